# Update countries & provincias Spain
# Refresh COVID-19 per-country figures and the "last updated" timestamp.
# A handful of countries changed rank (their row order swaps) once the
# new totals are in, so both the label (col A) and the stats move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 14:19"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 2356841
$ws.Range("C4").Value = 184
$ws.Range("E4").Value = 1254224
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 122250

# Alemania (row 14)
$ws.Range("B14").Value = 191653
$ws.Range("C14").Value = 78
$ws.Range("E14").Value = 7791

# Bielorrusia (row 27)
$ws.Range("B27").Value = 59023
$ws.Range("C27").Value = 518
$ws.Range("D27").Value = 37923
$ws.Range("E27").Value = 20749
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 351

# Paises Bajos (row 31)
$ws.Range("B31").Value = 49658
$ws.Range("C31").Value = 65

# Filipinas (row 43)
$ws.Range("B43").Value = 30682
$ws.Range("C43").Value = 630
$ws.Range("D43").Value = 8143
$ws.Range("E43").Value = 21362
$ws.Range("G43").Value = 8
$ws.Range("H43").Value = 1177

# Barein (row 50)
$ws.Range("E50").Value = 5281
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 64

# Corea del Sur / Dinamarca swap ranks (rows 63-64)
$ws.Range("A63").Value = "Dinamarca"
$ws.Range("B63").Value = 12527
$ws.Range("C63").Value = 136
$ws.Range("D63").Value = 11347
$ws.Range("E63").Value = 578
$ws.Range("G63").Value = 2
$ws.Range("H63").Value = 602

$ws.Range("A64").Value = "Corea del Sur"
$ws.Range("B64").Value = 12438
$ws.Range("C64").Value = 17
$ws.Range("D64").Value = 10881
$ws.Range("E64").Value = 1277
$ws.Range("H64").Value = 280

# Nepal (row 69)
$ws.Range("D69").Value = 2148
$ws.Range("E69").Value = 7390

# Malasia / Sudan swap ranks (rows 71-72)
$ws.Range("A71").Value = "Sudan"
$ws.Range("B71").Value = 8698
$ws.Range("C71").Value = 118
$ws.Range("D71").Value = 3460
$ws.Range("E71").Value = 4705
$ws.Range("G71").Value = 12
$ws.Range("H71").Value = 533

$ws.Range("A72").Value = "Malasia"
$ws.Range("B72").Value = 8587
$ws.Range("C72").Value = 15
$ws.Range("D72").Value = 8177
$ws.Range("E72").Value = 289
$ws.Range("H72").Value = 121

# Finlandia (row 75)
$ws.Range("D75").Value = 6400
$ws.Range("E75").Value = 417
$ws.Range("G75").Value = 1
$ws.Range("H75").Value = 327

# Croacia (row 101)
$ws.Range("B101").Value = 2336
$ws.Range("C101").Value = 19
$ws.Range("E101").Value = 87

# Mongolia / Siria swap ranks (rows 164-165)
$ws.Range("A164").Value = "Siria"
$ws.Range("B164").Value = 219
$ws.Range("C164").Value = 15
$ws.Range("D164").Value = 83
$ws.Range("E164").Value = 129
$ws.Range("H164").Value = 7

$ws.Range("A165").Value = "Mongolia"
$ws.Range("B165").Value = 213
$ws.Range("C165").Value = 7
$ws.Range("D165").Value = 153
$ws.Range("E165").Value = 60
$ws.Range("H165").Value = 0

# Dominica / Fiyi swap ranks (rows 202-203, stats stay tied at 18)
$ws.Range("A202").Value = "Fiyi"
$ws.Range("A203").Value = "Dominica"

# Montserrat / Seychelles swap ranks (rows 211-212)
$ws.Range("A211").Value = "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0

$ws.Range("A212").Value = "Montserrat"
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1

# Islas Virgenes Britanicas / Papua Nueva Guinea swap ranks (rows 214-215)
$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0

$ws.Range("A215").Value = "Islas Virgenes Britanicas"
$ws.Range("D215").Value = 7
$ws.Range("H215").Value = 1
